$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column F: "Minimap" / "string" / "小地图" headers, mirroring the
#     existing B/C/D header rows (row1 = English, row2 = type, row3 = Chinese).
#     Copy formatting from the matching B-column header cells first so the
#     new header cells land on the same existing shared style (s="1")
#     instead of minting a new one. ---
$ws.Range("B1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Minimap"

$ws.Range("B2").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("F2").Value = "string"

$ws.Range("B3").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$ws.Range("F3").Value = "小地图"

# --- Formulas for the minimap texture filename, mirroring column C.
#     Copy formatting from column C's data cells (style s="3") first so the
#     formula cells land on the same existing style, then fill in the
#     formula. F5:F6 are set together so they share one formula
#     definition, just like C5:C6 already do. ---
$ws.Range("C4").Copy()
$ws.Range("F4").PasteSpecial(-4122)
$ws.Range("F4").Formula = '=CONCATENATE("tex_habitats_minimap_",B4,".png")'

$ws.Range("C5:C6").Copy()
$ws.Range("F5:F6").PasteSpecial(-4122)
$ws.Range("F5:F6").Formula = '=CONCATENATE("tex_habitats_minimap_",B5,".png")'

$excel.CutCopyMode = 0

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 12.5
$ws.Columns.Item(2).ColumnWidth = 12.5
$ws.Columns.Item(3).ColumnWidth = 20.375
$ws.Columns.Item(4).ColumnWidth = 12.5
$ws.Columns.Item(5).ColumnWidth = 12.5
$ws.Columns.Item(6).ColumnWidth = 29.125

# --- Selection moves to B3 ---
$ws.Range("B3").Select()
